# Append the new daily log row (2025/10/08) as row 78, matching the
# existing data rows: date + weekday as text, hour + ranking as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 78

$dateCell = $ws.Cells.Item($row, 1)
# Force the date-like string to be stored as plain text (not an Excel
# date serial) so it round-trips the same way as the existing rows.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/08"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "水"
$ws.Cells.Item($row, 3).Value = 14
$ws.Cells.Item($row, 4).Value = 13
